$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-06-13"

# Update the header label in I1 for the running total column
$ws.Range("I1").Value = "2022 (through 06-13)"

# Update June (row 7) value for the 2022 running-total column
$ws.Range("I7").Value = 59

# Update Total row (row 14) value for the 2022 running-total column
$ws.Range("I14").Value = 722
